# Fruta / hortaliza, semanal
# Insert two new weekly rows at the top of the data block (row 19),
# pushing the existing rows down by two, then populate the two new rows
# with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 19..30 down to 21..32 by inserting two new blank rows at 19:20
$ws.Range("19:20").Insert()

# New row 19
$ws.Range("A19").Value() = 9
$ws.Range("B19").Value() = "Vega Central Mapocho de Santiago"
$ws.Range("C19").Value() = "Metropolitana"
$ws.Range("D19").Value() = 44435
$ws.Range("E19").Value() = 13
$ws.Range("F19").Value() = "Fruta"
$ws.Range("G19").Value() = 100108
$ws.Range("H19").Value() = "Tropicales y subtropicales"
$ws.Range("I19").Value() = 100108007
$ws.Range("J19").Value() = "Coco"
$ws.Range("K19").Value() = "Sin especificar"
$ws.Range("L19").Value() = "Primera"
$ws.Range("M19").Value() = 60
$ws.Range("N19").Value() = 25000
$ws.Range("O19").Value() = 25000
$ws.Range("P19").Value() = 25000
$ws.Range("Q19").Value() = "$/malla 20 unidades"
$ws.Range("R19").Value() = "Perú"
$ws.Range("S19").Value() = 1250
$ws.Range("T19").Value() = 20

# New row 20
$ws.Range("A20").Value() = 9
$ws.Range("B20").Value() = "Vega Central Mapocho de Santiago"
$ws.Range("C20").Value() = "Metropolitana"
$ws.Range("D20").Value() = 44431
$ws.Range("E20").Value() = 13
$ws.Range("F20").Value() = "Fruta"
$ws.Range("G20").Value() = 100108
$ws.Range("H20").Value() = "Tropicales y subtropicales"
$ws.Range("I20").Value() = 100108007
$ws.Range("J20").Value() = "Coco"
$ws.Range("K20").Value() = "Sin especificar"
$ws.Range("L20").Value() = "Primera"
$ws.Range("M20").Value() = 60
$ws.Range("N20").Value() = 25000
$ws.Range("O20").Value() = 25000
$ws.Range("P20").Value() = 25000
$ws.Range("Q20").Value() = "$/malla 20 unidades"
$ws.Range("R20").Value() = "Perú"
$ws.Range("S20").Value() = 1250
$ws.Range("T20").Value() = 20
